$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Planilha1")
$ws.Range("C4").Value = "communication failure; absence of procedure; inadequate supervision; insufficient training; lack of planning; high turnover; lack of resources; outdated procedure; weak security culture; ineffective management of change; failure in risk assessment; inappropriate prioritization; lack of lessons learned; technical review missing; production pressure; documentation inconsistency; no procedure; lack of communication; missing procedure; no schedule; poor maintenance; no inspection; unplanned activity; preventive maintenance issue; planning failure; no pre-task briefing; poor training"
$ws.Range("C4").Select()
